$d = $word.ActiveDocument

# Locate the "Lots of other tidbits on github" paragraph. The new DevNot
# paragraph needs to be inserted immediately before it.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Lots of other tidbits on*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find 'Lots of other tidbits on' paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)
$insertionPoint = $target.Range
$insertionPoint.Collapse(1)

# Insert a new (empty) paragraph before the target one. It inherits the
# "Definition" paragraph style used by the surrounding project bullets.
$insertionPoint.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($targetIndex)
$newRange = $newPara.Range
$newRange.Text = "DevNot: developer notes with code friendly markdown editor. Rails, Backbone."

$paraStart = $newPara.Range.Start
$linkRange = $d.Range($paraStart, $paraStart + 6)

$d.Hyperlinks.Add($linkRange, "https://devnot.com") | Out-Null

# The collection is ordered by document position, so re-scan for the
# hyperlink we just created (by address) and apply the same look used by
# the other project links: bold text styled with the "Link" character
# style.
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $hl = $d.Hyperlinks.Item($i)
    if ($hl.Address -eq "https://devnot.com/") {
        $hl.Range.Style = "Link"
        $hl.Range.Font.Bold = 1
        break
    }
}
